$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Step 1: populate B (set), C (species), D (count) cells for rows 216-294.
# Rows are written in the exact order the source data was entered so that new
# shared-string table entries land at the same indices as the target workbook
# (a couple of rows -- 228, and the 252/253 pair -- were filled out of strict
# row order in the original edit).
$ws.Range("B216").Value = "134"
$ws.Range("C216").Value = "coho"
$ws.Range("D216").Value = 1
$ws.Range("B217").Value = "134"
$ws.Range("C217").Value = "lingcod"
$ws.Range("D217").Value = 5
$ws.Range("B218").Value = "134"
$ws.Range("C218").Value = "canary rockfish"
$ws.Range("D218").Value = 1
$ws.Range("B219").Value = "135"
$ws.Range("C219").Value = "halibut"
$ws.Range("D219").Value = 1
$ws.Range("B220").Value = "135"
$ws.Range("C220").Value = "coho"
$ws.Range("D220").Value = 1
$ws.Range("B221").Value = "136"
$ws.Range("C221").Value = "lingcod"
$ws.Range("D221").Value = 1
$ws.Range("B222").Value = "136"
$ws.Range("C222").Value = "canary rockfish"
$ws.Range("D222").Value = 1
$ws.Range("B223").Value = "136"
$ws.Range("C223").Value = "coho"
$ws.Range("D223").Value = 1
$ws.Range("B224").Value = "137"
$ws.Range("C224").Value = "yellowtail rockfish"
$ws.Range("D224").Value = 2
$ws.Range("B225").Value = "137"
$ws.Range("C225").Value = "coho"
$ws.Range("D225").Value = 1
$ws.Range("B226").Value = "138"
$ws.Range("C226").Value = "coho"
$ws.Range("D226").Value = 1
$ws.Range("B227").Value = "139"
$ws.Range("C227").Value = "coho"
$ws.Range("D227").Value = 3
$ws.Range("B229").Value = "141"
$ws.Range("C229").Value = "coho"
$ws.Range("D229").Value = 1
$ws.Range("B230").Value = "142"
$ws.Range("C230").Value = "pink"
$ws.Range("D230").Value = 1
$ws.Range("B231").Value = "144"
$ws.Range("C231").Value = "coho"
$ws.Range("D231").Value = 1
$ws.Range("B232").Value = "145"
$ws.Range("C232").Value = "coho"
$ws.Range("D232").Value = 2
$ws.Range("B233").Value = "146"
$ws.Range("C233").Value = "coho"
$ws.Range("D233").Value = 2
$ws.Range("B234").Value = "147"
$ws.Range("C234").Value = "coho"
$ws.Range("D234").Value = 5
$ws.Range("B235").Value = "147"
$ws.Range("C235").Value = "pink"
$ws.Range("D235").Value = 1
$ws.Range("B236").Value = "148"
$ws.Range("C236").Value = "pink"
$ws.Range("D236").Value = 1
$ws.Range("B237").Value = "148"
$ws.Range("C237").Value = "coho"
$ws.Range("D237").Value = 5
$ws.Range("B238").Value = "149"
$ws.Range("C238").Value = "coho"
$ws.Range("D238").Value = 7
$ws.Range("B239").Value = "149"
$ws.Range("C239").Value = "pink"
$ws.Range("D239").Value = 1
$ws.Range("B240").Value = "149"
$ws.Range("C240").Value = "chinook"
$ws.Range("D240").Value = 2
$ws.Range("B241").Value = "150"
$ws.Range("C241").Value = "coho"
$ws.Range("D241").Value = 6
$ws.Range("B242").Value = "150"
$ws.Range("C242").Value = "pink"
$ws.Range("D242").Value = 1
$ws.Range("B243").Value = "150"
$ws.Range("C243").Value = "chinook"
$ws.Range("D243").Value = 1
$ws.Range("B244").Value = "151"
$ws.Range("C244").Value = "coho"
$ws.Range("D244").Value = 1
$ws.Range("B228").Value = "140"
$ws.Range("C228").Value = "chinook"
$ws.Range("D228").Value = 1
$ws.Range("B245").Value = "152"
$ws.Range("C245").Value = "chinook"
$ws.Range("D245").Value = 1
$ws.Range("B246").Value = "152"
$ws.Range("C246").Value = "pink"
$ws.Range("D246").Value = 1
$ws.Range("B247").Value = "153"
$ws.Range("C247").Value = "pink"
$ws.Range("D247").Value = 1
$ws.Range("B248").Value = "153"
$ws.Range("C248").Value = "coho"
$ws.Range("D248").Value = 1
$ws.Range("B249").Value = "153"
$ws.Range("C249").Value = "ling cod"
$ws.Range("D249").Value = 1
$ws.Range("B250").Value = "153"
$ws.Range("C250").Value = "sanddab"
$ws.Range("D250").Value = 1
$ws.Range("B251").Value = "154"
$ws.Range("C251").Value = "coho"
$ws.Range("D251").Value = 1
$ws.Range("B253").Value = "156"
$ws.Range("C253").Value = "pink"
$ws.Range("D253").Value = 1
$ws.Range("B252").Value = "155"
$ws.Range("C252").Value = "coho"
$ws.Range("D252").Value = 1
$ws.Range("B254").Value = "158"
$ws.Range("C254").Value = "chinook"
$ws.Range("D254").Value = 1
$ws.Range("B255").Value = "158"
$ws.Range("C255").Value = "lingcod"
$ws.Range("D255").Value = 1
$ws.Range("B256").Value = "158"
$ws.Range("C256").Value = "coho"
$ws.Range("D256").Value = 2
$ws.Range("B257").Value = "159"
$ws.Range("C257").Value = "coho"
$ws.Range("D257").Value = 5
$ws.Range("B258").Value = "159"
$ws.Range("C258").Value = "pink"
$ws.Range("D258").Value = 2
$ws.Range("B259").Value = "160"
$ws.Range("C259").Value = "coho"
$ws.Range("D259").Value = 3
$ws.Range("B260").Value = "161"
$ws.Range("C260").Value = "chinook"
$ws.Range("D260").Value = 1
$ws.Range("B261").Value = "161"
$ws.Range("C261").Value = "coho"
$ws.Range("D261").Value = 4
$ws.Range("B262").Value = "162"
$ws.Range("C262").Value = "coho"
$ws.Range("D262").Value = 3
$ws.Range("B263").Value = "163"
$ws.Range("C263").Value = "coho"
$ws.Range("D263").Value = 5
$ws.Range("B264").Value = "164"
$ws.Range("C264").Value = "coho"
$ws.Range("D264").Value = 3
$ws.Range("B265").Value = "166"
$ws.Range("C265").Value = "pink"
$ws.Range("D265").Value = 1
$ws.Range("B266").Value = "167"
$ws.Range("C266").Value = "coho"
$ws.Range("D266").Value = 1
$ws.Range("B267").Value = "168"
$ws.Range("C267").Value = "coho"
$ws.Range("D267").Value = 2
$ws.Range("B268").Value = "170"
$ws.Range("C268").Value = "chinook"
$ws.Range("D268").Value = 2
$ws.Range("B269").Value = "170"
$ws.Range("C269").Value = "coho"
$ws.Range("D269").Value = 2
$ws.Range("B270").Value = "171"
$ws.Range("C270").Value = "chinook"
$ws.Range("D270").Value = 1
$ws.Range("B271").Value = "171"
$ws.Range("C271").Value = "coho"
$ws.Range("D271").Value = 1
$ws.Range("B272").Value = "172"
$ws.Range("C272").Value = "coho"
$ws.Range("D272").Value = 1
$ws.Range("B273").Value = "173"
$ws.Range("C273").Value = "chinook"
$ws.Range("D273").Value = 1
$ws.Range("B274").Value = "173"
$ws.Range("C274").Value = "pink"
$ws.Range("D274").Value = 3
$ws.Range("B275").Value = "173"
$ws.Range("C275").Value = "coho"
$ws.Range("D275").Value = 2
$ws.Range("B276").Value = "174"
$ws.Range("C276").Value = "chinook"
$ws.Range("D276").Value = 2
$ws.Range("B277").Value = "174"
$ws.Range("C277").Value = "coho"
$ws.Range("D277").Value = 2
$ws.Range("B278").Value = "175"
$ws.Range("C278").Value = "chinook"
$ws.Range("D278").Value = 1
$ws.Range("B279").Value = "175"
$ws.Range("C279").Value = "coho"
$ws.Range("D279").Value = 7
$ws.Range("B280").Value = "176"
$ws.Range("C280").Value = "chinook"
$ws.Range("D280").Value = 1
$ws.Range("B281").Value = "176"
$ws.Range("C281").Value = "coho"
$ws.Range("D281").Value = 5
$ws.Range("B282").Value = "177"
$ws.Range("C282").Value = "coho"
$ws.Range("D282").Value = 3
$ws.Range("B283").Value = "178"
$ws.Range("C283").Value = "chinook"
$ws.Range("D283").Value = 2
$ws.Range("B284").Value = "178"
$ws.Range("C284").Value = "coho"
$ws.Range("D284").Value = 9
$ws.Range("B285").Value = "179"
$ws.Range("C285").Value = "chinook"
$ws.Range("D285").Value = 1
$ws.Range("B286").Value = "180"
$ws.Range("C286").Value = "coho"
$ws.Range("D286").Value = 8
$ws.Range("B287").Value = "182"
$ws.Range("C287").Value = "coho"
$ws.Range("D287").Value = 5
$ws.Range("B288").Value = "183"
$ws.Range("C288").Value = "chinook"
$ws.Range("D288").Value = 1
$ws.Range("B289").Value = "183"
$ws.Range("C289").Value = "coho"
$ws.Range("D289").Value = 3
$ws.Range("B290").Value = "184"
$ws.Range("C290").Value = "coho"
$ws.Range("D290").Value = 1
$ws.Range("B291").Value = "185"
$ws.Range("C291").Value = "chinook"
$ws.Range("D291").Value = 1
$ws.Range("B292").Value = "185"
$ws.Range("C292").Value = "coho"
$ws.Range("D292").Value = 7
$ws.Range("B293").Value = "186"
$ws.Range("C293").Value = "chinook"
$ws.Range("D293").Value = 2
$ws.Range("B294").Value = "186"
$ws.Range("C294").Value = "pink"
$ws.Range("D294").Value = 2

# Step 2: fill column A with the "CT2019"&"_"&B shared formula. Done as two
# range-level assignments (matching the two fill-down passes visible in the
# target workbook: one extending through row 242, a second for 243-294).
$ws.Range("A216:A242").Formula = '="CT2019"&"_"&B216'
$ws.Range("A243:A294").Formula = '="CT2019"&"_"&B243'

# Step 3: restore the view state (active cell) visible in the target workbook.
$ws.Range("C270").Select()
